$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- Fill values in the order that reproduces the shared-string table of the target file --

# New string #8: "Validte the functionality of the Login Account."
$ws.Range("D5").Value = "Validte the functionality of the Login Account."

# New string #9: "P0" (used by E4, E5, E6)
$ws.Range("E4").Value = "P0"

# New string #10: "(TS_002)`nLogin functionality"
$ws.Range("B5").Value = "(TS_002)" + [char]10 + "Login functionality"

# New string #11: "(TS_003)`nForget Password functionality"
$ws.Range("B6").Value = "(TS_003)" + [char]10 + "Forget Password functionality"

# New string #12: "Validte the functionality of the Forget Password."
$ws.Range("D6").Value = "Validte the functionality of the Forget Password."

# -- Fill remaining (already-existing string / numeric) values --
$ws.Range("C5").Value = "FRS"
$ws.Range("E5").Value = "P0"
$ws.Range("F5").Value = 6
$ws.Rows.Item(5).RowHeight = 30

$ws.Range("C6").Value = "FRS"
$ws.Range("E6").Value = "P0"
$ws.Range("F6").Value = 19
$ws.Rows.Item(6).RowHeight = 45

$ws.Range("F4").Value = 34

# Row 7 - Total formula
$ws.Range("F7").Formula = "=SUM(F4:F6)"

# Update selection to F7
$ws.Range("F7").Select()
